# New crime data collected
# Updates the weekly CompStat 63rd Precinct report: bumps the report
# volume/number and date range in the header, and refreshes the crime
# statistics table (rows 15-27) with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text (volume/number, week covering dates) ----
$ws.Range("A8").Value = "Volume 30   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/5/2023  Through  6/11/2023"

# ---- Cells that change numeric/text type need their style copied from ----
# ---- a donor cell that already carries the right number format,      ----
# ---- so we don't introduce brand-new cell styles.                    ----

# D16 & E16: were blank-like placeholders ("0" / "***.*" text), now real counts
$ws.Range("C16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 1

$ws.Range("H16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = 100

# C23: numeric count -> "0" placeholder text (style copied from C14, already "0")
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)

# D26: numeric count -> "0" placeholder text
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial(-4122)

# E26: numeric % change -> "***.*" placeholder text (non-numeric string, no
# extra trick needed, just re-point the style to the General/text style)
$ws.Range("C14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = "***.*"

# ---- Remaining plain value updates across the Week/28-Day/YTD/2-Year table ----
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 2
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 49
$ws.Range("K16").Value = -10.204081632653
$ws.Range("L16").Value = 46.666666666666
$ws.Range("M16").Value = -55.102040816326
$ws.Range("N16").Value = -85.084745762711
$ws.Range("F17").Value = 13
$ws.Range("H17").Value = 8.333333333333
$ws.Range("I17").Value = 71
$ws.Range("J17").Value = 71
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5.970149253731
$ws.Range("M17").Value = 24.561403508771
$ws.Range("N17").Value = -47.407407407407
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 28.571428571428
$ws.Range("I18").Value = 46
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = 2.222222222222
$ws.Range("L18").Value = 4.545454545454
$ws.Range("M18").Value = -57.407407407407
$ws.Range("N18").Value = -91.771019677996
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -8.333333333333
$ws.Range("I19").Value = 265
$ws.Range("J19").Value = 208
$ws.Range("K19").Value = 27.403846153846
$ws.Range("L19").Value = 102.290076335878
$ws.Range("M19").Value = 21.559633027522
$ws.Range("N19").Value = -12.541254125412
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = 46
$ws.Range("J20").Value = 51
$ws.Range("K20").Value = -9.803921568627
$ws.Range("L20").Value = 24.324324324324
$ws.Range("M20").Value = -36.986301369863
$ws.Range("N20").Value = -96.394984326018
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 16.666666666666
$ws.Range("F21").Value = 87
$ws.Range("H21").Value = 7.407407407407
$ws.Range("I21").Value = 478
$ws.Range("J21").Value = 428
$ws.Range("K21").Value = 11.682242990654
$ws.Range("L21").Value = 52.229299363057
$ws.Range("M21").Value = -14.490161001788
$ws.Range("N21").Value = -81.544401544401
$ws.Range("L23").Value = 40
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 18.181818181818
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = 12.5
$ws.Range("I24").Value = 528
$ws.Range("J24").Value = 464
$ws.Range("K24").Value = 13.793103448275
$ws.Range("L24").Value = 54.838709677419
$ws.Range("M24").Value = 30.69306930693
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -83.333333333333
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = -4.761904761904
$ws.Range("I25").Value = 114
$ws.Range("J25").Value = 96
$ws.Range("K25").Value = 18.75
$ws.Range("L25").Value = 4.587155963302
$ws.Range("M25").Value = -20.833333333333
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("L26").Value = 37.5
$ws.Range("D27").Value = 1
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -60
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = 13.333333333333
$ws.Range("L27").Value = 142.857142857143
